$d = $word.ActiveDocument
$find = $d.Content.Find

$find.Execute("33×69=2277", $true, $false, $false, $false, $false, $true, 1, $false, "50×51=2550", 2) | Out-Null
$find.Execute("22×61=1342", $true, $false, $false, $false, $false, $true, 1, $false, "23×47=1081", 2) | Out-Null
$find.Execute("96×39=3744", $true, $false, $false, $false, $false, $true, 1, $false, "68×91=6188", 2) | Out-Null
$find.Execute("85×20=1700", $true, $false, $false, $false, $false, $true, 1, $false, "100×93=9300", 2) | Out-Null
$find.Execute("63×16=1008", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=4761", 2) | Out-Null
$find.Execute("86×22=1892", $true, $false, $false, $false, $false, $true, 1, $false, "99×21=2079", 2) | Out-Null
$find.Execute("31×32=992", $true, $false, $false, $false, $false, $true, 1, $false, "76×26=1976", 2) | Out-Null
$find.Execute("80×64=5120", $true, $false, $false, $false, $false, $true, 1, $false, "13×65=845", 2) | Out-Null
$find.Execute("27×95=2565", $true, $false, $false, $false, $false, $true, 1, $false, "52×56=2912", 2) | Out-Null
$find.Execute("76×22=1672", $true, $false, $false, $false, $false, $true, 1, $false, "63×22=1386", 2) | Out-Null
$find.Execute("58×75=4350", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=7728", 2) | Out-Null
$find.Execute("77×84=6468", $true, $false, $false, $false, $false, $true, 1, $false, "41×50=2050", 2) | Out-Null
$find.Execute("28×23=644", $true, $false, $false, $false, $false, $true, 1, $false, "72×93=6696", 2) | Out-Null
$find.Execute("85×72=6120", $true, $false, $false, $false, $false, $true, 1, $false, "39×31=1209", 2) | Out-Null
$find.Execute("53×60=3180", $true, $false, $false, $false, $false, $true, 1, $false, "55×71=3905", 2) | Out-Null
$find.Execute("74×58=4292", $true, $false, $false, $false, $false, $true, 1, $false, "92×25=2300", 2) | Out-Null
$find.Execute("40×94=3760", $true, $false, $false, $false, $false, $true, 1, $false, "49×33=1617", 2) | Out-Null
$find.Execute("20×36=720", $true, $false, $false, $false, $false, $true, 1, $false, "11×48=528", 2) | Out-Null
$find.Execute("74×13=962", $true, $false, $false, $false, $false, $true, 1, $false, "100×39=3900", 2) | Out-Null
$find.Execute("70×99=6930", $true, $false, $false, $false, $false, $true, 1, $false, "44×79=3476", 2) | Out-Null
$find.Execute("48×30=1440", $true, $false, $false, $false, $false, $true, 1, $false, "47×64=3008", 2) | Out-Null
$find.Execute("17×45=765", $true, $false, $false, $false, $false, $true, 1, $false, "44×49=2156", 2) | Out-Null
$find.Execute("74×82=6068", $true, $false, $false, $false, $false, $true, 1, $false, "66×49=3234", 2) | Out-Null
$find.Execute("70×11=770", $true, $false, $false, $false, $false, $true, 1, $false, "74×66=4884", 2) | Out-Null
$find.Execute("72×59=4248", $true, $false, $false, $false, $false, $true, 1, $false, "99×69=6831", 2) | Out-Null
$find.Execute("11×33=363", $true, $false, $false, $false, $false, $true, 1, $false, "68×73=4964", 2) | Out-Null
$find.Execute("52×58=3016", $true, $false, $false, $false, $false, $true, 1, $false, "21×59=1239", 2) | Out-Null
$find.Execute("76×57=4332", $true, $false, $false, $false, $false, $true, 1, $false, "91×40=3640", 2) | Out-Null
$find.Execute("82×76=6232", $true, $false, $false, $false, $false, $true, 1, $false, "60×65=3900", 2) | Out-Null
$find.Execute("39×62=2418", $true, $false, $false, $false, $false, $true, 1, $false, "11×34=374", 2) | Out-Null
$find.Execute("72×61=4392", $true, $false, $false, $false, $false, $true, 1, $false, "51×61=3111", 2) | Out-Null
$find.Execute("71×30=2130", $true, $false, $false, $false, $false, $true, 1, $false, "36×73=2628", 2) | Out-Null
$find.Execute("34×65=2210", $true, $false, $false, $false, $false, $true, 1, $false, "27×43=1161", 2) | Out-Null
$find.Execute("13×77=1001", $true, $false, $false, $false, $false, $true, 1, $false, "31×86=2666", 2) | Out-Null
$find.Execute("49×80=3920", $true, $false, $false, $false, $false, $true, 1, $false, "41×58=2378", 2) | Out-Null
$find.Execute("37×54=1998", $true, $false, $false, $false, $false, $true, 1, $false, "85×38=3230", 2) | Out-Null
$find.Execute("50×26=1300", $true, $false, $false, $false, $false, $true, 1, $false, "49×53=2597", 2) | Out-Null
$find.Execute("45×61=2745", $true, $false, $false, $false, $false, $true, 1, $false, "72×81=5832", 2) | Out-Null
$find.Execute("96×54=5184", $true, $false, $false, $false, $false, $true, 1, $false, "34×84=2856", 2) | Out-Null
$find.Execute("68×86=5848", $true, $false, $false, $false, $false, $true, 1, $false, "23×46=1058", 2) | Out-Null
$find.Execute("13×18=234", $true, $false, $false, $false, $false, $true, 1, $false, "81×40=3240", 2) | Out-Null
$find.Execute("24×83=1992", $true, $false, $false, $false, $false, $true, 1, $false, "49×50=2450", 2) | Out-Null
$find.Execute("71×81=5751", $true, $false, $false, $false, $false, $true, 1, $false, "33×13=429", 2) | Out-Null
$find.Execute("75×69=5175", $true, $false, $false, $false, $false, $true, 1, $false, "70×10=700", 2) | Out-Null
$find.Execute("55×69=3795", $true, $false, $false, $false, $false, $true, 1, $false, "10×52=520", 2) | Out-Null
$find.Execute("46×22=1012", $true, $false, $false, $false, $false, $true, 1, $false, "55×18=990", 2) | Out-Null
$find.Execute("29×16=464", $true, $false, $false, $false, $false, $true, 1, $false, "91×37=3367", 2) | Out-Null
$find.Execute("23×60=1380", $true, $false, $false, $false, $false, $true, 1, $false, "57×60=3420", 2) | Out-Null
$find.Execute("42×11=462", $true, $false, $false, $false, $false, $true, 1, $false, "58×62=3596", 2) | Out-Null
$find.Execute("21×14=294", $true, $false, $false, $false, $false, $true, 1, $false, "43×99=4257", 2) | Out-Null
$find.Execute("40×73=2920", $true, $false, $false, $false, $false, $true, 1, $false, "68×66=4488", 2) | Out-Null
$find.Execute("41×24=984", $true, $false, $false, $false, $false, $true, 1, $false, "98×61=5978", 2) | Out-Null
$find.Execute("93×94=8742", $true, $false, $false, $false, $false, $true, 1, $false, "27×68=1836", 2) | Out-Null
$find.Execute("94×50=4700", $true, $false, $false, $false, $false, $true, 1, $false, "29×74=2146", 2) | Out-Null
$find.Execute("42×30=1260", $true, $false, $false, $false, $false, $true, 1, $false, "42×31=1302", 2) | Out-Null
$find.Execute("45×51=2295", $true, $false, $false, $false, $false, $true, 1, $false, "76×43=3268", 2) | Out-Null
$find.Execute("53×19=1007", $true, $false, $false, $false, $false, $true, 1, $false, "27×85=2295", 2) | Out-Null
$find.Execute("86×67=5762", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=6014", 2) | Out-Null
$find.Execute("68×50=3400", $true, $false, $false, $false, $false, $true, 1, $false, "51×84=4284", 2) | Out-Null
$find.Execute("88×43=3784", $true, $false, $false, $false, $false, $true, 1, $false, "81×34=2754", 2) | Out-Null
$find.Execute("41×28=1148", $true, $false, $false, $false, $false, $true, 1, $false, "88×33=2904", 2) | Out-Null
$find.Execute("18×13=234", $true, $false, $false, $false, $false, $true, 1, $false, "11×19=209", 2) | Out-Null
$find.Execute("67×93=6231", $true, $false, $false, $false, $false, $true, 1, $false, "69×63=4347", 2) | Out-Null
$find.Execute("75×93=6975", $true, $false, $false, $false, $false, $true, 1, $false, "35×75=2625", 2) | Out-Null
$find.Execute("13×28=364", $true, $false, $false, $false, $false, $true, 1, $false, "50×46=2300", 2) | Out-Null
$find.Execute("26×93=2418", $true, $false, $false, $false, $false, $true, 1, $false, "46×20=920", 2) | Out-Null
$find.Execute("42×73=3066", $true, $false, $false, $false, $false, $true, 1, $false, "91×100=9100", 2) | Out-Null
$find.Execute("11×30=330", $true, $false, $false, $false, $false, $true, 1, $false, "16×41=656", 2) | Out-Null
$find.Execute("28×34=952", $true, $false, $false, $false, $false, $true, 1, $false, "33×25=825", 2) | Out-Null
$find.Execute("98×83=8134", $true, $false, $false, $false, $false, $true, 1, $false, "99×92=9108", 2) | Out-Null
$find.Execute("30×68=2040", $true, $false, $false, $false, $false, $true, 1, $false, "49×45=2205", 2) | Out-Null
$find.Execute("20×93=1860", $true, $false, $false, $false, $false, $true, 1, $false, "58×49=2842", 2) | Out-Null
$find.Execute("59×77=4543", $true, $false, $false, $false, $false, $true, 1, $false, "34×88=2992", 2) | Out-Null
$find.Execute("44×46=2024", $true, $false, $false, $false, $false, $true, 1, $false, "41×91=3731", 2) | Out-Null
$find.Execute("32×11=352", $true, $false, $false, $false, $false, $true, 1, $false, "40×19=760", 2) | Out-Null
$find.Execute("34×32=1088", $true, $false, $false, $false, $false, $true, 1, $false, "98×100=9800", 2) | Out-Null
$find.Execute("83×50=4150", $true, $false, $false, $false, $false, $true, 1, $false, "29×53=1537", 2) | Out-Null
$find.Execute("44×92=4048", $true, $false, $false, $false, $false, $true, 1, $false, "97×65=6305", 2) | Out-Null
$find.Execute("22×74=1628", $true, $false, $false, $false, $false, $true, 1, $false, "67×90=6030", 2) | Out-Null
$find.Execute("16×71=1136", $true, $false, $false, $false, $false, $true, 1, $false, "12×59=708", 2) | Out-Null
$find.Execute("92×85=7820", $true, $false, $false, $false, $false, $true, 1, $false, "49×53=2597", 2) | Out-Null
$find.Execute("45×36=1620", $true, $false, $false, $false, $false, $true, 1, $false, "17×100=1700", 2) | Out-Null
$find.Execute("79×60=4740", $true, $false, $false, $false, $false, $true, 1, $false, "42×77=3234", 2) | Out-Null
$find.Execute("35×19=665", $true, $false, $false, $false, $false, $true, 1, $false, "20×62=1240", 2) | Out-Null
$find.Execute("84×26=2184", $true, $false, $false, $false, $false, $true, 1, $false, "42×55=2310", 2) | Out-Null
$find.Execute("10×34=340", $true, $false, $false, $false, $false, $true, 1, $false, "25×62=1550", 2) | Out-Null
$find.Execute("12×16=192", $true, $false, $false, $false, $false, $true, 1, $false, "86×34=2924", 2) | Out-Null
$find.Execute("10×95=950", $true, $false, $false, $false, $false, $true, 1, $false, "18×10=180", 2) | Out-Null
$find.Execute("53×79=4187", $true, $false, $false, $false, $false, $true, 1, $false, "44×80=3520", 2) | Out-Null
$find.Execute("81×24=1944", $true, $false, $false, $false, $false, $true, 1, $false, "26×50=1300", 2) | Out-Null
$find.Execute("44×43=1892", $true, $false, $false, $false, $false, $true, 1, $false, "21×79=1659", 2) | Out-Null
$find.Execute("45×66=2970", $true, $false, $false, $false, $false, $true, 1, $false, "34×40=1360", 2) | Out-Null
$find.Execute("36×93=3348", $true, $false, $false, $false, $false, $true, 1, $false, "97×30=2910", 2) | Out-Null
$find.Execute("50×28=1400", $true, $false, $false, $false, $false, $true, 1, $false, "29×36=1044", 2) | Out-Null
$find.Execute("96×98=9408", $true, $false, $false, $false, $false, $true, 1, $false, "25×100=2500", 2) | Out-Null
$find.Execute("41×89=3649", $true, $false, $false, $false, $false, $true, 1, $false, "45×12=540", 2) | Out-Null
$find.Execute("79×96=7584", $true, $false, $false, $false, $false, $true, 1, $false, "47×41=1927", 2) | Out-Null
$find.Execute("10×100=1000", $true, $false, $false, $false, $false, $true, 1, $false, "95×46=4370", 2) | Out-Null
$find.Execute("63×69=4347", $true, $false, $false, $false, $false, $true, 1, $false, "65×82=5330", 2) | Out-Null
$find.Execute("54×65=3510", $true, $false, $false, $false, $false, $true, 1, $false, "59×30=1770", 2) | Out-Null
